$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added for "Macroferia Regional de Talca -
# Zanahoria" at the top of its (date-descending) block, which starts at row
# 489. Insert a fresh row there; Excel pushes rows 489:513 down to 490:514
# and grows the sheet's used range to A1:R514.
$ws.Rows.Item(489).Insert()

$ws.Cells.Item(489, 1).Value = 5
$ws.Cells.Item(489, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(489, 3).Value = "Maule"
$ws.Cells.Item(489, 4).Value = 45041
$ws.Cells.Item(489, 5).Value = 7
$ws.Cells.Item(489, 6).Value = 100114013
$ws.Cells.Item(489, 7).Value = "Zanahoria"
$ws.Cells.Item(489, 8).Value = "Sin especificar"
$ws.Cells.Item(489, 9).Value = "Primera"
$ws.Cells.Item(489, 10).Value = 500
$ws.Cells.Item(489, 11).Value = 6000
$ws.Cells.Item(489, 12).Value = 6000
$ws.Cells.Item(489, 13).Value = 6000
$ws.Cells.Item(489, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(489, 15).Value = "Región de Ñuble"
$ws.Cells.Item(489, 16).Value = 300
$ws.Cells.Item(489, 17).Value = 20
$ws.Cells.Item(489, 18).Value = "Hortaliza"

# Match the date-cell number formatting used by the rest of column D.
$ws.Cells.Item(489, 4).NumberFormat = $ws.Cells.Item(490, 4).NumberFormat
